$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30 (the existing 1997 row), shifting rows 30-53 down to 31-54
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new data point
$ws.Cells.Item(30, 1).Value = 1997
$ws.Cells.Item(30, 2).Value = 8284431.7
